$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(29)
$full = $p.Range
$firstChar = $d.Range($full.Start, $full.Start+1)
$firstChar.Font.Color = 255
Write-Output "colored"
$insertPos = $d.Range($full.End - 2, $full.End - 2)
$insertPos.InsertAfter("X")
Write-Output "after insert"
